$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45895
$ws.Range("B2").Value = 5341.91564027273
$ws.Range("C2").Value = 4396.3097727477
$ws.Range("D2").Value = 4728
$ws.Range("E2").Value = 5689.108108
$ws.Range("F2").Value = 0.645926686456808

$ws.Range("A3").Value = 45896
$ws.Range("B3").Value = 5341.91530752175
$ws.Range("C3").Value = 4373.14886167392
$ws.Range("D3").Value = 1944
$ws.Range("E3").Value = 5689.107764
$ws.Range("F3").Value = 115.68088825634

$ws.Range("A4").Value = 45897
$ws.Range("B4").Value = 5341.91604070457
$ws.Range("C4").Value = 4363.67126786433
$ws.Range("D4").Value = 1944
$ws.Range("E4").Value = 5689.108523
$ws.Range("F4").Value = 115.28598958999

$ws.Range("A5").Value = 45898
$ws.Range("B5").Value = 5341.91604070659
$ws.Range("C5").Value = 3709.11736259734
$ws.Range("D5").Value = 1944
$ws.Range("E5").Value = 5689.108523
$ws.Range("F5").Value = 88.0129102037811

$ws.Range("A6").Value = 45899
$ws.Range("B6").Value = 975.205470885634
$ws.Range("C6").Value = 1587.55485165417
$ws.Range("D6").Value = 1944
$ws.Range("E6").Value = 1881.077314
$ws.Range("F6").Value = 22.892778948689

$ws.Range("A7").Value = 45900
$ws.Range("B7").Value = 832.128236149678
$ws.Range("C7").Value = 1544.3359306026
$ws.Range("D7").Value = 1944
$ws.Range("E7").Value = 1682.350742
$ws.Range("F7").Value = 18.7732681855383

$ws.Range("A8").Value = 45901
$ws.Range("B8").Value = 5776.55841866516
$ws.Range("C8").Value = 4969.12759581159
$ws.Range("D8").Value = 2952
$ws.Range("E8").Value = 6391.166845
$ws.Range("F8").Value = 109.655667589435

$ws.Range("A9").Value = 45902
$ws.Range("B9").Value = 5784.11073458408
$ws.Range("C9").Value = 5160.29477647573
$ws.Range("D9").Value = 2952
$ws.Range("E9").Value = 6398.783217
$ws.Range("F9").Value = 117.623635787152

$ws.Range("A10").Value = 45903
$ws.Range("B10").Value = 5784.11073458408
$ws.Range("C10").Value = 5136.67407143287
$ws.Range("D10").Value = 2952
$ws.Range("E10").Value = 6398.783217
$ws.Range("F10").Value = 116.6394397437

$ws.Range("A11").Value = 45904
$ws.Range("B11").Value = 5812.79253171584
$ws.Range("C11").Value = 5164.8804111703
$ws.Range("D11").Value = 2952
$ws.Range("E11").Value = 6445.475115
$ws.Range("F11").Value = 118.565124768936

$ws.Range("A12").Value = 45905
$ws.Range("B12").Value = 5812.79253171584
$ws.Range("C12").Value = 4463.36487502072
$ws.Range("D12").Value = 2952
$ws.Range("E12").Value = 6445.475115
$ws.Range("F12").Value = 89.3353107627034

$ws.Range("A13").Value = 45906
$ws.Range("B13").Value = 1096.28225605742
$ws.Range("C13").Value = 2196.52829077272
$ws.Range("D13").Value = 2952
$ws.Range("E13").Value = 2113.135789
$ws.Range("F13").Value = 10.8909093214708

$ws.Range("A14").Value = 45907
$ws.Range("B14").Value = 970.493609902267
$ws.Range("C14").Value = 2141.36193841852
$ws.Range("D14").Value = 2952
$ws.Range("E14").Value = 1974.395268
$ws.Range("F14").Value = 8.05264985484388

$ws.Range("A15").Value = 45908
$ws.Range("B15").Value = 6110.60514061168
$ws.Range("C15").Value = 5034.18791521718
$ws.Range("D15").Value = 2952
$ws.Range("E15").Value = 6751.192496
$ws.Range("F15").Value = 113.448969608563

